# Update cosinor analysis results for CircadiPy simulations (sine_0.1)
# Re-run of CircaDB / CircadiPy analyses with refreshed output values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.91000000000061
$ws.Range("H2").Value = 0.0007820139022649553
$ws.Range("I2").Value = 0.0007820139022649553
$ws.Range("L2").Value = 41.0773617164706
$ws.Range("M2").Value = '[17.401006639046997, 64.7537167938942]'
$ws.Range("N2").Value = 0.0010791853051908
$ws.Range("O2").Value = 0.0010791853051908
$ws.Range("P2").Value = 1.515763422452733
$ws.Range("Q2").Value = '[0.735868549489501, 2.295658295415965]'
$ws.Range("R2").Value = 0.0003040243337464332
$ws.Range("S2").Value = 0.0003040243337464332
$ws.Range("T2").Value = 54.04821233448838
$ws.Range("U2").Value = '[39.29663702121307, 68.79978764776368]'
$ws.Range("V2").Value = "2.798421272487417e-09"
$ws.Range("W2").Value = "2.798421272487417e-09"
$ws.Range("X2").Value = 19.6594394394399
$ws.Range("Y2").Value = 16.44338338338377
$ws.Range("Z2").Value = 22.87549549549604
$ws.Range("F3").Value = 25.91000000000061
$ws.Range("H3").Value = 0.005752086396564549
$ws.Range("I3").Value = 0.005752086396564549
$ws.Range("L3").Value = 34.86870132241917
$ws.Range("M3").Value = '[7.81899796221397, 61.91840468262437]'
$ws.Range("N3").Value = 0.01268385083173507
$ws.Range("O3").Value = 0.01268385083173507
$ws.Range("P3").Value = 1.163552834662887
$ws.Range("Q3").Value = '[0.30818426431611545, 2.0189214050096576]'
$ws.Range("R3").Value = 0.008782567571823607
$ws.Range("S3").Value = 0.008782567571823607
$ws.Range("T3").Value = 63.7969829525307
$ws.Range("U3").Value = '[49.071017756519794, 78.5229481485416]'
$ws.Range("V3").Value = "3.10456105268031e-11"
$ws.Range("W3").Value = "3.10456105268031e-11"
$ws.Range("X3").Value = 21.11185185185235
$ws.Range("Y3").Value = 17.58456456456498
$ws.Range("Z3").Value = 24.63913913913972
$ws.Range("F4").Value = 25.91000000000061
$ws.Range("H4").Value = 0.02733594760966118
$ws.Range("I4").Value = 0.02733594760966118
$ws.Range("L4").Value = 23.39429861915409
$ws.Range("M4").Value = '[1.0008402853600558, 45.78775695294813]'
$ws.Range("N4").Value = 0.04098706308458389
$ws.Range("O4").Value = 0.04098706308458389
$ws.Range("P4").Value = 1.213868632918579
$ws.Range("Q4").Value = '[-0.1195000208572683, 2.5472372866944273]'
$ws.Range("R4").Value = 0.07333323350246945
$ws.Range("S4").Value = 0.07333323350246945
$ws.Range("T4").Value = 53.42181968690402
$ws.Range("U4").Value = '[41.32505998935385, 65.51857938445417]'
$ws.Range("V4").Value = "1.785549486044147e-11"
$ws.Range("W4").Value = "1.785549486044147e-11"
$ws.Range("X4").Value = 20.90436436436485
$ws.Range("Y4").Value = 15.40594594594631
$ws.Range("Z4").Value = 26.4027827827834
$ws.Range("F5").Value = 25.91000000000061
$ws.Range("H5").Value = 0.006662222926533623
$ws.Range("I5").Value = 0.006662222926533623
$ws.Range("L5").Value = 36.17825276064104
$ws.Range("M5").Value = '[9.401021556532363, 62.955483964749725]'
$ws.Range("N5").Value = 0.00921572154375716
$ws.Range("O5").Value = 0.00921572154375716
$ws.Range("P5").Value = -0.1257894956392311
$ws.Range("Q5").Value = '[-1.0440528138056173, 0.792473822527155]'
$ws.Range("R5").Value = 0.7838844890985095
$ws.Range("S5").Value = 0.7838844890985095
$ws.Range("T5").Value = 53.22762413457428
$ws.Range("U5").Value = '[38.24669357881305, 68.2085546903355]'
$ws.Range("V5").Value = "5.981769524154856e-09"
$ws.Range("W5").Value = "5.981769524154856e-09"
$ws.Range("X5").Value = 0.5187187187187305
$ws.Range("Y5").Value = -3.267927927928006
$ws.Range("Z5").Value = 4.305365365365467
$ws.Range("B6").Value = 0
$ws.Range("F6").Value = 25.91000000000061
$ws.Range("H6").Value = 0.1494641433816739
$ws.Range("I6").Value = 0.1494641433816739
$ws.Range("L6").Value = 23.20338036670267
$ws.Range("M6").Value = '[-4.360287888647498, 50.76704862205284]'
$ws.Range("N6").Value = 0.09689147898662642
$ws.Range("O6").Value = 0.09689147898662642
$ws.Range("P6").Value = -0.01257894956392303
$ws.Range("Q6").Value = '[-1.616395018964118, 1.591237119836272]'
$ws.Range("R6").Value = 0.9874662616607934
$ws.Range("S6").Value = 0.9874662616607934
$ws.Range("T6").Value = 66.28589228571191
$ws.Range("U6").Value = '[50.22671118692497, 82.34507338449885]'
$ws.Range("V6").Value = "1.211244438081849e-10"
$ws.Range("W6").Value = "1.211244438081849e-10"
$ws.Range("X6").Value = 0.05187187187187092
$ws.Range("Y6").Value = -6.561791791791951
$ws.Range("Z6").Value = 6.665535535535692
$ws.Range("B7").Value = 1
$ws.Range("F7").Value = 25.91000000000061
$ws.Range("H7").Value = "7.563083680439675e-07"
$ws.Range("I7").Value = "7.563083680439675e-07"
$ws.Range("L7").Value = 54.40170665746451
$ws.Range("M7").Value = '[30.09308756881427, 78.71032574611475]'
$ws.Range("N7").Value = "4.64602129783831e-05"
$ws.Range("O7").Value = "4.64602129783831e-05"
$ws.Range("P7").Value = -0.7170001251436169
$ws.Range("Q7").Value = '[-1.1446844103170015, -0.28931583997023225]'
$ws.Range("R7").Value = 0.001521660505534594
$ws.Range("S7").Value = 0.001521660505534594
$ws.Range("T7").Value = 59.57867777511596
$ws.Range("U7").Value = '[46.95473729268238, 72.20261825754955]'
$ws.Range("V7").Value = "2.484457084506175e-12"
$ws.Range("W7").Value = "2.484457084506175e-12"
$ws.Range("X7").Value = 2.95669669669677
$ws.Range("Y7").Value = 1.193053053053088
$ws.Range("Z7").Value = 4.720340340340453
$ws.Range("F8").Value = 25.91000000000061
$ws.Range("H8").Value = "4.344789919674685e-05"
$ws.Range("I8").Value = "4.344789919674685e-05"
$ws.Range("L8").Value = 41.68353717081867
$ws.Range("M8").Value = '[19.721732440826557, 63.64534190081079]'
$ws.Range("N8").Value = 0.0004029615002092424
$ws.Range("O8").Value = 0.0004029615002092424
$ws.Range("P8").Value = -1.232737057264464
$ws.Range("Q8").Value = '[-1.8239476867688493, -0.6415264277600778]'
$ws.Range("R8").Value = 0.0001246363550244656
$ws.Range("S8").Value = 0.0001246363550244656
$ws.Range("T8").Value = 54.91676134383549
$ws.Range("U8").Value = '[42.59429439011565, 67.23922829755533]'
$ws.Range("V8").Value = "1.369415691954146e-11"
$ws.Range("W8").Value = "1.369415691954146e-11"
$ws.Range("X8").Value = 5.083443443443564
$ws.Range("Y8").Value = 2.645465465465527
$ws.Range("Z8").Value = 7.5214214214216
$ws.Range("F9").Value = 23.70000000000027
$ws.Range("H9").Value = 0.01224896922734453
$ws.Range("I9").Value = 0.01224896922734453
$ws.Range("L9").Value = 35.40337015207322
$ws.Range("M9").Value = '[6.610130536461057, 64.19660976768539]'
$ws.Range("N9").Value = 0.01709416049090784
$ws.Range("O9").Value = 0.01709416049090784
$ws.Range("P9").Value = -1.761052938949233
$ws.Range("Q9").Value = '[-2.8428426014466197, -0.6792632764518469]'
$ws.Range("R9").Value = 0.002015509661379422
$ws.Range("S9").Value = 0.002015509661379422
$ws.Range("T9").Value = 73.99000680886905
$ws.Range("U9").Value = '[57.91883253406456, 90.06118108367355]'
$ws.Range("V9").Value = "5.242473122279989e-12"
$ws.Range("W9").Value = "5.242473122279989e-12"
$ws.Range("X9").Value = 6.642642642642716
$ws.Range("Y9").Value = 2.562162162162189
$ws.Range("Z9").Value = 10.72312312312324
$ws.Range("F10").Value = 23.70000000000027
$ws.Range("H10").Value = 0.01097950806201076
$ws.Range("I10").Value = 0.01097950806201076
$ws.Range("L10").Value = 28.39534786848823
$ws.Range("M10").Value = '[5.408960868900017, 51.38173486807645]'
$ws.Range("N10").Value = 0.01661483099970207
$ws.Range("O10").Value = 0.01661483099970207
$ws.Range("P10").Value = -1.270473905956233
$ws.Range("Q10").Value = '[-2.276789871070081, -0.2641579408423844]'
$ws.Range("R10").Value = 0.01450573925398402
$ws.Range("S10").Value = 0.01450573925398402
$ws.Range("T10").Value = 53.00946977660055
$ws.Range("U10").Value = '[40.28933178337502, 65.72960776982609]'
$ws.Range("V10").Value = "9.285239244150034e-11"
$ws.Range("W10").Value = "9.285239244150034e-11"
$ws.Range("X10").Value = 4.792192192192246
$ws.Range("Y10").Value = 0.9963963963964053
$ws.Range("Z10").Value = 8.587987987988086
$ws.Range("F11").Value = 23.70000000000027
$ws.Range("H11").Value = 0.01100957816451686
$ws.Range("I11").Value = 0.01100957816451686
$ws.Range("L11").Value = 34.52115276245787
$ws.Range("M11").Value = '[8.285055604297732, 60.757249920618015]'
$ws.Range("N11").Value = 0.01106401991802541
$ws.Range("O11").Value = 0.01106401991802541
$ws.Range("P11").Value = -1.673000292001772
$ws.Range("Q11").Value = '[-2.7296320553713125, -0.6163685286322318]'
$ws.Range("R11").Value = 0.002599451648451323
$ws.Range("S11").Value = 0.002599451648451323
$ws.Range("T11").Value = 68.74213471362916
$ws.Range("U11").Value = '[53.277514222182525, 84.20675520507578]'
$ws.Range("V11").Value = "1.476818667356383e-11"
$ws.Range("W11").Value = "1.476818667356383e-11"
$ws.Range("X11").Value = 6.310510510510582
$ws.Range("Y11").Value = 2.324924924924951
$ws.Range("Z11").Value = 10.29609609609621
